$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Cort-Sstr4 -> ECs)
$ws.Range("M2").Value = 6.365199
$ws.Range("N2").Value = 19.095597
$ws.Range("O2").Value = 0.8551243431489016
$ws.Range("P2").Value = 0.8551243431489015
$ws.Range("Q2").Value = 0.700265246252
$ws.Range("R2").Value = 6.302387216268
$ws.Range("S2").Value = 0.8551243431489016
$ws.Range("T2").Value = 0.8551243431489015

# Row 3 (Cort-Sstr4 -> FAPs)
$ws.Range("O3").Value = 0.06239073379600107
$ws.Range("P3").Value = 0.06239073379600107
$ws.Range("S3").Value = 0.06239073379600107
$ws.Range("T3").Value = 0.06239073379600107

# Row 4 (Cort-Sstr4 -> MuSCs)
$ws.Range("M4").Value = 0.4859933333333333
$ws.Range("N4").Value = 1.45798
$ws.Range("O4").Value = 0.06529013938785132
$ws.Range("P4").Value = 0.06529013938785132
$ws.Range("Q4").Value = 0.05346639456888889
$ws.Range("R4").Value = 0.48119755112
$ws.Range("S4").Value = 0.06529013938785132
$ws.Range("T4").Value = 0.06529013938785132

# Row 5 (Cort-Sstr4 -> Resolving-Mac)
$ws.Range("M5").Value = 0.127991
$ws.Range("N5").Value = 0.383973
$ws.Range("O5").Value = 0.01719478366724608
$ws.Range("P5").Value = 0.01719478366724608
$ws.Range("Q5").Value = 0.01408088720133333
$ws.Range("R5").Value = 0.126727984812
$ws.Range("S5").Value = 0.01719478366724608
$ws.Range("T5").Value = 0.01719478366724608

$wb.Save()
